# Insert two new data rows (weekly update) right before the existing
# row that starts the 2021-11-xx block, shifting all subsequent rows down
# by 2 (so former row 520 becomes 522, ..., former row 622 becomes 624).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 520-521; everything from old row 520
# downward shifts to 522 downward (Excel's native row-insert semantics).
$ws.Rows("520:521").Insert()

# --- New row 520 : Primera, fecha 44694 ---
$ws.Range("A520").Value = 9
$ws.Range("B520").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C520").Value = "Metropolitana"
$ws.Range("D520").Value = 44694
$ws.Range("E520").Value = 13
$ws.Range("F520").Value = 100112008
$ws.Range("G520").Value = "Coliflor"
$ws.Range("H520").Value = "Sin especificar"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 2500
$ws.Range("K520").Value = 900
$ws.Range("L520").Value = 1000
$ws.Range("M520").Value = 950
$ws.Range("N520").Value = "`$/unidad"
$ws.Range("O520").Value = "Región Metropolitana"
$ws.Range("P520").Value = 950
$ws.Range("Q520").Value = 1
$ws.Range("R520").Value = "Hortaliza"

# --- New row 521 : Segunda, fecha 44694 ---
$ws.Range("A521").Value = 9
$ws.Range("B521").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C521").Value = "Metropolitana"
$ws.Range("D521").Value = 44694
$ws.Range("E521").Value = 13
$ws.Range("F521").Value = 100112008
$ws.Range("G521").Value = "Coliflor"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Segunda"
$ws.Range("J521").Value = 1060
$ws.Range("K521").Value = 800
$ws.Range("L521").Value = 800
$ws.Range("M521").Value = 800
$ws.Range("N521").Value = "`$/unidad"
$ws.Range("O521").Value = "Región Metropolitana"
$ws.Range("P521").Value = 800
$ws.Range("Q521").Value = 1
$ws.Range("R521").Value = "Hortaliza"
